# Gamepad_AT_2022_11142022.pptx — bug fixes & tuned sleep time for open/close
# claw & Bump Y function: add the missing "Left Bumper, Right Bumper & B ..."
# instruction line under the Tophat navigation-mode bullet on slide 1, and
# grow the text box to fit the new line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the instructions text box (Google Shape;122;p1) by its persistent
# shape Id rather than a hard-coded collection index.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 122) {
        $targetShape = $candidate
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the last paragraph ("Left Trigger, Right Trigger & A to put robot
# Tophat in navigation mode") and append a brand-new paragraph after it.
# NB: Paragraphs().Count (no explicit start arg) reports the true count;
# Paragraphs(1,-1).Count under-counts by one in this host, so use the
# no-arg form just to size the loop / find the last index.
$paraCount = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($paraCount, 1)

$newText = "Left Bumper, Right Bumper & B is to turn robot left into 180 degrees"
$null = $lastPara.InsertAfter("`r" + $newText)

# Re-fetch the text range/paragraph so formatting edits actually stick
# (edits applied to the range object returned by InsertAfter do not persist).
$tr2 = $targetShape.TextFrame.TextRange
$newParaCount = $tr2.Paragraphs().Count
$newPara = $tr2.Paragraphs($newParaCount, 1)

# Run 1: "Left Bumper, Right Bumper & B is to " -> blue, explicit plain style
$run1 = $newPara.Characters(1, 36)
$run1.Font.Bold = 0
$run1.Font.Italic = 0
$run1.Font.Underline = 0
$run1.Font.Strikethrough = 0
$run1.Font.Caps = 0
$run1.Font.Color.RGB = 15773696   # 00B0F0

# Run 2: "turn robo" -> red, explicit plain style
$run2 = $newPara.Characters(37, 9)
$run2.Font.Bold = 0
$run2.Font.Italic = 0
$run2.Font.Underline = 0
$run2.Font.Strikethrough = 0
$run2.Font.Caps = 0
$run2.Font.Color.RGB = 255        # FF0000

# Run 3: "t left " -> red
$run3 = $newPara.Characters(46, 7)
$run3.Font.Color.RGB = 255        # FF0000

# Run 4: "into 180 degrees" -> blue
$run4 = $newPara.Characters(53, 16)
$run4.Font.Color.RGB = 15773696   # 00B0F0

# Grow the (auto-fit) text box to accommodate the extra line.
$targetShape.Height = 156.308545
